$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
# (Overview!E2/F2 mirror the same status string and pick it up automatically,
# but we set them explicitly too so every copy of the text is refreshed.)
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"

# --- Latest Handback DateTime refreshed for the new handback ---
$wsZhCn.Range("K2").Value = "2016-08-21 18:53:12"
$wsDeDe.Range("K2").Value = "2016-08-21 18:53:19"

# --- Error Detail cleared: handback is now in sync, no stale-version error ---
$wsZhCn.Range("P2").Value = ""
$wsDeDe.Range("P2").Value = ""

# --- Column widths widened to fit the longer "Handed back: in sync with en-US" status text ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.17
$wsOverview.Columns.Item(6).ColumnWidth = 29.17
$wsZhCn.Columns.Item(3).ColumnWidth = 29.17
$wsDeDe.Columns.Item(3).ColumnWidth = 29.17

# --- Error Detail column narrowed now that it is empty ---
$wsZhCn.Columns.Item(16).ColumnWidth = 12.8
$wsDeDe.Columns.Item(16).ColumnWidth = 12.8
